$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '254.75'
    'E2' = '3.61%'
    'D3' = '28.27'
    'E3' = '-6.41%'
    'D4' = '5.247'
    'E4' = '1.79%'
    'D5' = '0.05878'
    'E5' = '1.99%'
    'D6' = '6.722'
    'E6' = '0.73%'
    'B7' = 'GateToken'
    'C7' = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
    'D7' = '3.223'
    'E7' = '-0.81%'
    'B8' = 'MXToken'
    'C8' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'D8' = '0.8657'
    'E8' = '1.95%'
    'B9' = 'FTXToken'
    'C9' = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
    'D9' = '0.9883'
    'E9' = '15.22%'
    'B10' = 'WazirX'
    'C10' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'D10' = '0.1409'
    'E10' = '1.90%'
    'B11' = 'MandalaExchangeToken'
    'C11' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'D11' = '0.07171'
    'E11' = '1.19%'
    'B12' = 'BitrueCoin'
    'C12' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'D12' = '0.03188'
    'E12' = '-2.31%'
    'B13' = 'BitMartToken'
    'C13' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'D13' = '0.09227'
    'E13' = '-1.47%'
    'B14' = 'BitForexToken'
    'C14' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'D14' = '0.001547'
    'E14' = '1.19%'
    'B15' = 'One'
    'C15' = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
    'D15' = '0.0006094'
    'E15' = '2.65%'
    'B16' = 'TigerCash'
    'C16' = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    'D16' = '0.005826'
    'E16' = '-3.33%'
    'B17' = 'LEO'
    'C17' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'D17' = '3.499'
    'E17' = '-0.77%'
    'E18' = '-0.58%'
    'D19' = '0.3180'
    'E19' = '1.75%'
    'E20' = '2.35%'
    'D21' = '0.1293'
    'E21' = '-1.68%'
    'D22' = '3.542'
    'E22' = '1.30%'
    'D23' = '0.04157'
    'E23' = '0.80%'
    'D25' = '0.001227'
    'E25' = '0.18%'
    'D26' = '0.004797'
    'E26' = '15.57%'
    'D27' = '0.0001201'
    'E27' = '0.09%'
    'D28' = '0.0001467'
    'E28' = '1.30%'
    'E40' = '1.54%'
    'B41' = 'BKEXToken'
    'C41' = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
    'D41' = '0.1101'
    'E41' = '2.79%'
    'B42' = 'KickToken'
    'C42' = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
    'D42' = '0.003808'
    'E42' = '-33.74%'
    'D43' = '0.002340'
    'E43' = '-5.24%'
    'D44' = '0.009702'
    'D45' = '0.00005233'
    'E45' = '-4.56%'
    'D46' = '0.00000000751'
    'E46' = '0.12%'
    'D47' = '0.1201'
    'E47' = '69.19%'
    'D48' = '0.002145'
    'E48' = '-13.06%'
    'E49' = '0.12%'
    'E50' = '0.12%'
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cellRef]
}
